$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.535.67"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "3.593.79"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "236.03"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "654.70"
$ws.Range("E6").Value = "  +4.75%  "
$ws.Range("D7").Value = "1.46"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("D8").Value = "0.401"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").Value = "3.591.39"
$ws.Range("E11").Value = "  +3.10%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "42.38"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "6.46"
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").Value = "4.285.00"
$ws.Range("E15").Value = "  +3.64%  "
$ws.Range("D16").Value = "95.342.60"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "0.0000253"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "3.596.63"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("D19").Value = "7.91"
$ws.Range("E19").Value = "  -5.00%  "
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "17.91"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "3.52"
$ws.Range("E22").Value = "  +4.01%  "
$ws.Range("D23").Value = "507.55"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").Value = "0.477"
$ws.Range("E24").Value = "  -5.06%  "
$ws.Range("D25").Value = "0.0000195"
$ws.Range("E25").Value = "  +6.24%  "
$ws.Range("D26").Value = "6.60"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "95.26"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("D28").Value = "3.788.67"
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").Value = "12.54"
$ws.Range("E29").Value = "  +2.80%  "
$ws.Range("D30").Value = "3.05"
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "11.25"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "0.139"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").Value = "32.37"
$ws.Range("E35").Value = "  +9.05%  "
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").Value = "0.560"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "8.12"
$ws.Range("E38").Value = "  +7.84%  "
$ws.Range("D39").Value = "557.46"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").Value = "35.65"
$ws.Range("E44").Value = "  +39.81%  "
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("D46").Value = "2.30"
$ws.Range("E46").Value = "  +6.59%  "
$ws.Range("D47").Value = "23.57"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").Value = "5.66"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("D49").Value = "0.0412"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").Value = "3.57"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").Value = "53.24"
$ws.Range("E51").Value = "  -0.03%  "
